$d = $word.ActiveDocument

$d.Content.Find.Execute("741÷6=123, 3", $true, $false, $false, $false, $false, $true, 1, $false, "373÷6=62, 1", 2)
$d.Content.Find.Execute("145÷6=24, 1", $true, $false, $false, $false, $false, $true, 1, $false, "540÷5=108, 0", 2)
$d.Content.Find.Execute("744÷4=186, 0", $true, $false, $false, $false, $false, $true, 1, $false, "149÷2=74, 1", 2)
$d.Content.Find.Execute("281÷4=70, 1", $true, $false, $false, $false, $false, $true, 1, $false, "745÷6=124, 1", 2)
$d.Content.Find.Execute("205÷2=102, 1", $true, $false, $false, $false, $false, $true, 1, $false, "894÷2=447, 0", 2)
$d.Content.Find.Execute("324÷5=64, 4", $true, $false, $false, $false, $false, $true, 1, $false, "455÷4=113, 3", 2)
$d.Content.Find.Execute("973÷5=194, 3", $true, $false, $false, $false, $false, $true, 1, $false, "579÷8=72, 3", 2)
$d.Content.Find.Execute("434÷7=62, 0", $true, $false, $false, $false, $false, $true, 1, $false, "373÷3=124, 1", 2)
$d.Content.Find.Execute("294÷9=32, 6", $true, $false, $false, $false, $false, $true, 1, $false, "614÷6=102, 2", 2)
$d.Content.Find.Execute("394÷3=131, 1", $true, $false, $false, $false, $false, $true, 1, $false, "310÷9=34, 4", 2)
$d.Content.Find.Execute("372÷6=62, 0", $true, $false, $false, $false, $false, $true, 1, $false, "992÷6=165, 2", 2)
$d.Content.Find.Execute("768÷7=109, 5", $true, $false, $false, $false, $false, $true, 1, $false, "710÷3=236, 2", 2)
$d.Content.Find.Execute("477÷5=95, 2", $true, $false, $false, $false, $false, $true, 1, $false, "600÷4=150, 0", 2)
$d.Content.Find.Execute("362÷8=45, 2", $true, $false, $false, $false, $false, $true, 1, $false, "906÷6=151, 0", 2)
$d.Content.Find.Execute("304÷5=60, 4", $true, $false, $false, $false, $false, $true, 1, $false, "145÷7=20, 5", 2)
$d.Content.Find.Execute("866÷4=216, 2", $true, $false, $false, $false, $false, $true, 1, $false, "898÷2=449, 0", 2)
$d.Content.Find.Execute("423÷6=70, 3", $true, $false, $false, $false, $false, $true, 1, $false, "209÷7=29, 6", 2)
$d.Content.Find.Execute("798÷2=399, 0", $true, $false, $false, $false, $false, $true, 1, $false, "571÷4=142, 3", 2)
$d.Content.Find.Execute("138÷6=23, 0", $true, $false, $false, $false, $false, $true, 1, $false, "598÷3=199, 1", 2)
$d.Content.Find.Execute("489÷6=81, 3", $true, $false, $false, $false, $false, $true, 1, $false, "208÷6=34, 4", 2)
$d.Content.Find.Execute("597÷8=74, 5", $true, $false, $false, $false, $false, $true, 1, $false, "800÷6=133, 2", 2)
$d.Content.Find.Execute("119÷3=39, 2", $true, $false, $false, $false, $false, $true, 1, $false, "270÷9=30, 0", 2)
$d.Content.Find.Execute("851÷5=170, 1", $true, $false, $false, $false, $false, $true, 1, $false, "391÷5=78, 1", 2)
$d.Content.Find.Execute("224÷2=112, 0", $true, $false, $false, $false, $false, $true, 1, $false, "501÷2=250, 1", 2)
$d.Content.Find.Execute("887÷4=221, 3", $true, $false, $false, $false, $false, $true, 1, $false, "911÷2=455, 1", 2)
